$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1025
$ws1.Range("F3").Value = 293
$ws1.Range("F4").Value = 1416
$ws1.Range("F5").Value = 8537
$ws1.Range("F6").Value = 67
$ws1.Range("F7").Value = 480
$ws1.Range("F9").Value = 252
$ws1.Range("F11").Value = 3436
$ws1.Range("F13").Value = 346
$ws1.Range("F14").Value = 63
$ws1.Range("F15").Value = 991
$ws1.Range("F17").Value = 1094
$ws1.Range("F18").Value = 296
$ws1.Range("F19").Value = 167
$ws1.Range("F20").Value = 2116

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 33

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1025
$ws4.Range("F3").Value = 293
$ws4.Range("F4").Value = 1416
$ws4.Range("F5").Value = 8537
$ws4.Range("F6").Value = 67
$ws4.Range("F7").Value = 480
$ws4.Range("F9").Value = 252
$ws4.Range("F11").Value = 3436
$ws4.Range("F13").Value = 346
$ws4.Range("F14").Value = 63
$ws4.Range("F15").Value = 991
$ws4.Range("F17").Value = 1094
$ws4.Range("F18").Value = 296
$ws4.Range("F19").Value = 167
$ws4.Range("F20").Value = 2116
$ws4.Range("F21").Value = 33
